# Refresh the cryptos price/volume table (GitHub Actions daily snapshot).
# Cells that look like plain decimal numbers ("226.09", "0.557", ...) are
# forced back to text via NumberFormat "@" before the write (otherwise COM
# auto-coerces the string into a float and trailing zeros / precision are
# lost), then the style is reset to "Normal" so no stray number-format is
# left attached to the cell - matches the original inline-string cells,
# which carry no explicit style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.421.44'
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").Value = '1.791.19'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.557'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.61%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.73'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.05%  '
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("D12").Value = '2.048.75'
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.10'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.58%  '
$ws.Range("D14").Value = '1.795.84'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.635'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.70%  '
$ws.Range("D16").Value = '34.380.59'
$ws.Range("E16").Value = '  +0.75%  '
$ws.Range("E17").Value = '  +1.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.73'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '247.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.30%  '
$ws.Range("E20").Value = '  +2.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.42%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("E23").Value = '  +1.11%  '
$ws.Range("E24").Value = '  +1.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '164.90'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.73%  '
$ws.Range("E27").Value = '  +0.93%  '
$ws.Range("E28").Value = '  +2.23%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.79'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.84%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("E32").Value = '  +0.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.90'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.27%  '
$ws.Range("E34").Value = '  +1.18%  '
$ws.Range("D35").Value = '1.418.42'
$ws.Range("E35").Value = '  -2.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.59'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.92%  '
$ws.Range("E37").Value = '  +2.65%  '
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("E39").Value = '  +1.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '84.79'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.40%  '
$ws.Range("E41").Value = '  +1.02%  '
$ws.Range("E42").Value = '  +1.21%  '
$ws.Range("E43").Value = '  +2.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.61'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.94%  '
$ws.Range("E45").Value = '  +2.91%  '
$ws.Range("E46").Value = '  -0.51%  '
$ws.Range("E47").Value = '  +0.26%  '
$ws.Range("D48").Value = '1.949.31'
$ws.Range("E48").Value = '  +0.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.22%  '
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0129'
$ws.Range("E51").Value = '  -4.87%  '
